# Apply the cell-value and view changes described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# --- Update boolean parameter cells in row 4 ---
$ws.Range("E4").Value = $true
$ws.Range("H4").Value = $false
$ws.Range("J4").Value = $false
$ws.Range("K4").Value = $false

# --- Update frozen-pane view / selection state ---
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E4:F4").Select()
